# Total Factor Productivity worksheet update:
# Add the missing 2018 data point at the top of the series (the sheet
# previously started at 2019), shifting the existing 2019-2023 rows down
# by one and extending the used range to A1:B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new row (row 7, currently blank) the same formatting as the
# existing data rows by copying row 6's format into it.
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)

# Shift the existing year/value pairs down one row (bottom-up so we don't
# clobber a value before it has been copied).
$ws.Range("A7").Value = $ws.Range("A6").Value2
$ws.Range("B7").Value = $ws.Range("B6").Value2
$ws.Range("A6").Value = $ws.Range("A5").Value2
$ws.Range("B6").Value = $ws.Range("B5").Value2
$ws.Range("A5").Value = $ws.Range("A4").Value2
$ws.Range("B5").Value = $ws.Range("B4").Value2
$ws.Range("A4").Value = $ws.Range("A3").Value2
$ws.Range("B4").Value = $ws.Range("B3").Value2
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2

# Insert the new 2018 data point into the now-vacated first data row.
$ws.Range("A2").Value = 2018
$ws.Range("B2").Value = 100.633

# Move the active cell/selection as recorded in the saved view state.
$null = $ws.Range("E14").Select()

# Refresh the "Generated on" timestamp in the footer.
$ws.PageSetup.RightFooter = "Generated on: November 5, 2024 (12:57:40 AM)"
